$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 437, shifting existing rows 437-559 down to 439-561
$ws.Rows(437).Resize(2).Insert()

# Populate the two brand-new rows (437, 438) completely
$ws.Cells.Item(437, 1).Value = 7
$ws.Cells.Item(437, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(437, 3).Value = "Ñuble"
$ws.Cells.Item(437, 4).Value = 44642
$ws.Cells.Item(437, 5).Value = 16
$ws.Cells.Item(437, 6).Value = "Fruta"
$ws.Cells.Item(437, 7).Value = 100108
$ws.Cells.Item(437, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(437, 9).Value = 100108006
$ws.Cells.Item(437, 10).Value = "Plátano"
$ws.Cells.Item(437, 11).Value = "Sin especificar"
$ws.Cells.Item(437, 12).Value = "Pintón"
$ws.Cells.Item(437, 13).Value = 200
$ws.Cells.Item(437, 14).Value = 19000
$ws.Cells.Item(437, 15).Value = 19000
$ws.Cells.Item(437, 16).Value = 19000
$ws.Cells.Item(437, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(437, 18).Value = "Ecuador"
$ws.Cells.Item(437, 19).Value = 950
$ws.Cells.Item(437, 20).Value = 20

$ws.Cells.Item(438, 1).Value = 7
$ws.Cells.Item(438, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(438, 3).Value = "Ñuble"
$ws.Cells.Item(438, 4).Value = 44642
$ws.Cells.Item(438, 5).Value = 16
$ws.Cells.Item(438, 6).Value = "Fruta"
$ws.Cells.Item(438, 7).Value = 100108
$ws.Cells.Item(438, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(438, 9).Value = 100108006
$ws.Cells.Item(438, 10).Value = "Plátano"
$ws.Cells.Item(438, 11).Value = "Sin especificar"
$ws.Cells.Item(438, 12).Value = "Primera Pintón"
$ws.Cells.Item(438, 13).Value = 300
$ws.Cells.Item(438, 14).Value = 20000
$ws.Cells.Item(438, 15).Value = 21000
$ws.Cells.Item(438, 16).Value = 20500
$ws.Cells.Item(438, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(438, 18).Value = "Ecuador"
$ws.Cells.Item(438, 19).Value = 1025
$ws.Cells.Item(438, 20).Value = 20

# Update D, L, M, N, O, P, S for the shifted rows (439-561) to their new target values
$ws.Cells.Item(439, 4).Value = 44473
$ws.Cells.Item(439, 12).Value = "Pintón"
$ws.Cells.Item(439, 13).Value = 150
$ws.Cells.Item(439, 14).Value = 13000
$ws.Cells.Item(439, 15).Value = 13000
$ws.Cells.Item(439, 16).Value = 13000
$ws.Cells.Item(439, 19).Value = 650
$ws.Cells.Item(440, 4).Value = 44473
$ws.Cells.Item(440, 12).Value = "Primera Pintón"
$ws.Cells.Item(440, 13).Value = 440
$ws.Cells.Item(440, 14).Value = 14000
$ws.Cells.Item(440, 15).Value = 15000
$ws.Cells.Item(440, 16).Value = 14500
$ws.Cells.Item(440, 19).Value = 725
$ws.Cells.Item(441, 4).Value = 44357
$ws.Cells.Item(441, 12).Value = "Pintón"
$ws.Cells.Item(441, 13).Value = 80
$ws.Cells.Item(441, 14).Value = 11000
$ws.Cells.Item(441, 15).Value = 11000
$ws.Cells.Item(441, 16).Value = 11000
$ws.Cells.Item(441, 19).Value = 550
$ws.Cells.Item(442, 4).Value = 44357
$ws.Cells.Item(442, 12).Value = "Primera Pintón"
$ws.Cells.Item(442, 13).Value = 160
$ws.Cells.Item(442, 14).Value = 12000
$ws.Cells.Item(442, 15).Value = 13000
$ws.Cells.Item(442, 16).Value = 12500
$ws.Cells.Item(442, 19).Value = 625
$ws.Cells.Item(443, 4).Value = 44455
$ws.Cells.Item(443, 12).Value = "Pintón"
$ws.Cells.Item(443, 13).Value = 300
$ws.Cells.Item(443, 14).Value = 18000
$ws.Cells.Item(443, 15).Value = 19000
$ws.Cells.Item(443, 16).Value = 18500
$ws.Cells.Item(443, 19).Value = 925
$ws.Cells.Item(444, 4).Value = 44455
$ws.Cells.Item(444, 12).Value = "Primera Pintón"
$ws.Cells.Item(444, 13).Value = 400
$ws.Cells.Item(444, 14).Value = 20000
$ws.Cells.Item(444, 15).Value = 21000
$ws.Cells.Item(444, 16).Value = 20500
$ws.Cells.Item(444, 19).Value = 1025
$ws.Cells.Item(445, 4).Value = 44581
$ws.Cells.Item(445, 12).Value = "Pintón"
$ws.Cells.Item(445, 13).Value = 80
$ws.Cells.Item(445, 14).Value = 14000
$ws.Cells.Item(445, 15).Value = 14000
$ws.Cells.Item(445, 16).Value = 14000
$ws.Cells.Item(445, 19).Value = 700
$ws.Cells.Item(446, 4).Value = 44581
$ws.Cells.Item(446, 12).Value = "Primera Pintón"
$ws.Cells.Item(446, 13).Value = 160
$ws.Cells.Item(446, 14).Value = 15000
$ws.Cells.Item(446, 15).Value = 16000
$ws.Cells.Item(446, 16).Value = 15500
$ws.Cells.Item(446, 19).Value = 775
$ws.Cells.Item(447, 4).Value = 44537
$ws.Cells.Item(447, 12).Value = "Pintón"
$ws.Cells.Item(447, 13).Value = 150
$ws.Cells.Item(447, 14).Value = 17000
$ws.Cells.Item(447, 15).Value = 17000
$ws.Cells.Item(447, 16).Value = 17000
$ws.Cells.Item(447, 19).Value = 850
$ws.Cells.Item(448, 4).Value = 44537
$ws.Cells.Item(448, 12).Value = "Primera Pintón"
$ws.Cells.Item(448, 13).Value = 500
$ws.Cells.Item(448, 14).Value = 18000
$ws.Cells.Item(448, 15).Value = 19000
$ws.Cells.Item(448, 16).Value = 18500
$ws.Cells.Item(448, 19).Value = 925
$ws.Cells.Item(449, 4).Value = 44553
$ws.Cells.Item(449, 12).Value = "Pintón"
$ws.Cells.Item(449, 13).Value = 200
$ws.Cells.Item(449, 14).Value = 11000
$ws.Cells.Item(449, 15).Value = 11000
$ws.Cells.Item(449, 16).Value = 11000
$ws.Cells.Item(449, 19).Value = 550
$ws.Cells.Item(450, 4).Value = 44553
$ws.Cells.Item(450, 12).Value = "Primera Pintón"
$ws.Cells.Item(450, 13).Value = 400
$ws.Cells.Item(450, 14).Value = 12000
$ws.Cells.Item(450, 15).Value = 13000
$ws.Cells.Item(450, 16).Value = 12500
$ws.Cells.Item(450, 19).Value = 625
$ws.Cells.Item(451, 4).Value = 44490
$ws.Cells.Item(451, 12).Value = "Pintón"
$ws.Cells.Item(451, 13).Value = 120
$ws.Cells.Item(451, 14).Value = 24000
$ws.Cells.Item(451, 15).Value = 24000
$ws.Cells.Item(451, 16).Value = 24000
$ws.Cells.Item(451, 19).Value = 1200
$ws.Cells.Item(452, 4).Value = 44490
$ws.Cells.Item(452, 12).Value = "Primera Pintón"
$ws.Cells.Item(452, 13).Value = 240
$ws.Cells.Item(452, 14).Value = 25000
$ws.Cells.Item(452, 15).Value = 26000
$ws.Cells.Item(452, 16).Value = 25500
$ws.Cells.Item(452, 19).Value = 1275
$ws.Cells.Item(453, 4).Value = 44397
$ws.Cells.Item(453, 12).Value = "Pintón"
$ws.Cells.Item(453, 13).Value = 80
$ws.Cells.Item(453, 14).Value = 11000
$ws.Cells.Item(453, 15).Value = 11000
$ws.Cells.Item(453, 16).Value = 11000
$ws.Cells.Item(453, 19).Value = 550
$ws.Cells.Item(454, 4).Value = 44397
$ws.Cells.Item(454, 12).Value = "Primera Pintón"
$ws.Cells.Item(454, 13).Value = 240
$ws.Cells.Item(454, 14).Value = 12000
$ws.Cells.Item(454, 15).Value = 13000
$ws.Cells.Item(454, 16).Value = 12500
$ws.Cells.Item(454, 19).Value = 625
$ws.Cells.Item(455, 4).Value = 44446
$ws.Cells.Item(455, 12).Value = "Pintón"
$ws.Cells.Item(455, 13).Value = 150
$ws.Cells.Item(455, 14).Value = 16000
$ws.Cells.Item(455, 15).Value = 16000
$ws.Cells.Item(455, 16).Value = 16000
$ws.Cells.Item(455, 19).Value = 800
$ws.Cells.Item(456, 4).Value = 44446
$ws.Cells.Item(456, 12).Value = "Primera Pintón"
$ws.Cells.Item(456, 13).Value = 240
$ws.Cells.Item(456, 14).Value = 17000
$ws.Cells.Item(456, 15).Value = 18000
$ws.Cells.Item(456, 16).Value = 17500
$ws.Cells.Item(456, 19).Value = 875
$ws.Cells.Item(457, 4).Value = 44641
$ws.Cells.Item(457, 12).Value = "Pintón"
$ws.Cells.Item(457, 13).Value = 320
$ws.Cells.Item(457, 14).Value = 19000
$ws.Cells.Item(457, 15).Value = 21000
$ws.Cells.Item(457, 16).Value = 20125
$ws.Cells.Item(457, 19).Value = 1006
$ws.Cells.Item(458, 4).Value = 44421
$ws.Cells.Item(458, 12).Value = "Pintón"
$ws.Cells.Item(458, 13).Value = 300
$ws.Cells.Item(458, 14).Value = 12000
$ws.Cells.Item(458, 15).Value = 12000
$ws.Cells.Item(458, 16).Value = 12000
$ws.Cells.Item(458, 19).Value = 600
$ws.Cells.Item(459, 4).Value = 44421
$ws.Cells.Item(459, 12).Value = "Primera Pintón"
$ws.Cells.Item(459, 13).Value = 600
$ws.Cells.Item(459, 14).Value = 13000
$ws.Cells.Item(459, 15).Value = 14000
$ws.Cells.Item(459, 16).Value = 13500
$ws.Cells.Item(459, 19).Value = 675
$ws.Cells.Item(460, 4).Value = 44329
$ws.Cells.Item(460, 12).Value = "Pintón"
$ws.Cells.Item(460, 13).Value = 120
$ws.Cells.Item(460, 14).Value = 14000
$ws.Cells.Item(460, 15).Value = 14000
$ws.Cells.Item(460, 16).Value = 14000
$ws.Cells.Item(460, 19).Value = 700
$ws.Cells.Item(461, 4).Value = 44329
$ws.Cells.Item(461, 12).Value = "Primera Pintón"
$ws.Cells.Item(461, 13).Value = 240
$ws.Cells.Item(461, 14).Value = 15000
$ws.Cells.Item(461, 15).Value = 16000
$ws.Cells.Item(461, 16).Value = 15500
$ws.Cells.Item(461, 19).Value = 775
$ws.Cells.Item(462, 4).Value = 44637
$ws.Cells.Item(462, 12).Value = "Pintón"
$ws.Cells.Item(462, 13).Value = 100
$ws.Cells.Item(462, 14).Value = 17000
$ws.Cells.Item(462, 15).Value = 17000
$ws.Cells.Item(462, 16).Value = 17000
$ws.Cells.Item(462, 19).Value = 850
$ws.Cells.Item(463, 4).Value = 44637
$ws.Cells.Item(463, 12).Value = "Primera Pintón"
$ws.Cells.Item(463, 13).Value = 300
$ws.Cells.Item(463, 14).Value = 18000
$ws.Cells.Item(463, 15).Value = 19000
$ws.Cells.Item(463, 16).Value = 18500
$ws.Cells.Item(463, 19).Value = 925
$ws.Cells.Item(464, 4).Value = 44208
$ws.Cells.Item(464, 12).Value = "Pintón"
$ws.Cells.Item(464, 13).Value = 840
$ws.Cells.Item(464, 14).Value = 12000
$ws.Cells.Item(464, 15).Value = 13000
$ws.Cells.Item(464, 16).Value = 12524
$ws.Cells.Item(464, 19).Value = 626
$ws.Cells.Item(465, 4).Value = 44208
$ws.Cells.Item(465, 12).Value = "Primera Pintón"
$ws.Cells.Item(465, 13).Value = 450
$ws.Cells.Item(465, 14).Value = 13500
$ws.Cells.Item(465, 15).Value = 14000
$ws.Cells.Item(465, 16).Value = 13778
$ws.Cells.Item(465, 19).Value = 689
$ws.Cells.Item(466, 4).Value = 44445
$ws.Cells.Item(466, 12).Value = "Pintón"
$ws.Cells.Item(466, 13).Value = 300
$ws.Cells.Item(466, 14).Value = 13500
$ws.Cells.Item(466, 15).Value = 14000
$ws.Cells.Item(466, 16).Value = 13750
$ws.Cells.Item(466, 19).Value = 688
$ws.Cells.Item(467, 4).Value = 44445
$ws.Cells.Item(467, 12).Value = "Primera Pintón"
$ws.Cells.Item(467, 13).Value = 300
$ws.Cells.Item(467, 14).Value = 14500
$ws.Cells.Item(467, 15).Value = 15000
$ws.Cells.Item(467, 16).Value = 14750
$ws.Cells.Item(467, 19).Value = 738
$ws.Cells.Item(468, 4).Value = 44524
$ws.Cells.Item(468, 12).Value = "Pintón"
$ws.Cells.Item(468, 13).Value = 80
$ws.Cells.Item(468, 14).Value = 20000
$ws.Cells.Item(468, 15).Value = 20000
$ws.Cells.Item(468, 16).Value = 20000
$ws.Cells.Item(468, 19).Value = 1000
$ws.Cells.Item(469, 4).Value = 44524
$ws.Cells.Item(469, 12).Value = "Primera Pintón"
$ws.Cells.Item(469, 13).Value = 120
$ws.Cells.Item(469, 14).Value = 21000
$ws.Cells.Item(469, 15).Value = 22000
$ws.Cells.Item(469, 16).Value = 21500
$ws.Cells.Item(469, 19).Value = 1075
$ws.Cells.Item(470, 4).Value = 44355
$ws.Cells.Item(470, 12).Value = "Pintón"
$ws.Cells.Item(470, 13).Value = 80
$ws.Cells.Item(470, 14).Value = 11000
$ws.Cells.Item(470, 15).Value = 11000
$ws.Cells.Item(470, 16).Value = 11000
$ws.Cells.Item(470, 19).Value = 550
$ws.Cells.Item(471, 4).Value = 44355
$ws.Cells.Item(471, 12).Value = "Primera Pintón"
$ws.Cells.Item(471, 13).Value = 160
$ws.Cells.Item(471, 14).Value = 12000
$ws.Cells.Item(471, 15).Value = 13000
$ws.Cells.Item(471, 16).Value = 12500
$ws.Cells.Item(471, 19).Value = 625
$ws.Cells.Item(472, 4).Value = 44530
$ws.Cells.Item(472, 12).Value = "Pintón"
$ws.Cells.Item(472, 13).Value = 160
$ws.Cells.Item(472, 14).Value = 18000
$ws.Cells.Item(472, 15).Value = 19000
$ws.Cells.Item(472, 16).Value = 18500
$ws.Cells.Item(472, 19).Value = 925
$ws.Cells.Item(473, 4).Value = 44530
$ws.Cells.Item(473, 12).Value = "Primera"
$ws.Cells.Item(473, 13).Value = 120
$ws.Cells.Item(473, 14).Value = 21000
$ws.Cells.Item(473, 15).Value = 21000
$ws.Cells.Item(473, 16).Value = 21000
$ws.Cells.Item(473, 19).Value = 1050
$ws.Cells.Item(474, 4).Value = 44530
$ws.Cells.Item(474, 12).Value = "Primera Pintón"
$ws.Cells.Item(474, 13).Value = 120
$ws.Cells.Item(474, 14).Value = 20000
$ws.Cells.Item(474, 15).Value = 20000
$ws.Cells.Item(474, 16).Value = 20000
$ws.Cells.Item(474, 19).Value = 1000
$ws.Cells.Item(475, 4).Value = 44483
$ws.Cells.Item(475, 12).Value = "Pintón"
$ws.Cells.Item(475, 13).Value = 300
$ws.Cells.Item(475, 14).Value = 19500
$ws.Cells.Item(475, 15).Value = 20000
$ws.Cells.Item(475, 16).Value = 19750
$ws.Cells.Item(475, 19).Value = 988
$ws.Cells.Item(476, 4).Value = 44483
$ws.Cells.Item(476, 12).Value = "Primera Pintón"
$ws.Cells.Item(476, 13).Value = 400
$ws.Cells.Item(476, 14).Value = 21000
$ws.Cells.Item(476, 15).Value = 22000
$ws.Cells.Item(476, 16).Value = 21500
$ws.Cells.Item(476, 19).Value = 1075
$ws.Cells.Item(477, 4).Value = 44294
$ws.Cells.Item(477, 12).Value = "Pintón"
$ws.Cells.Item(477, 13).Value = 80
$ws.Cells.Item(477, 14).Value = 13000
$ws.Cells.Item(477, 15).Value = 13000
$ws.Cells.Item(477, 16).Value = 13000
$ws.Cells.Item(477, 19).Value = 650
$ws.Cells.Item(478, 4).Value = 44294
$ws.Cells.Item(478, 12).Value = "Primera Pintón"
$ws.Cells.Item(478, 13).Value = 160
$ws.Cells.Item(478, 14).Value = 14000
$ws.Cells.Item(478, 15).Value = 15000
$ws.Cells.Item(478, 16).Value = 14500
$ws.Cells.Item(478, 19).Value = 725
$ws.Cells.Item(479, 4).Value = 44617
$ws.Cells.Item(479, 12).Value = "Pintón"
$ws.Cells.Item(479, 13).Value = 80
$ws.Cells.Item(479, 14).Value = 16000
$ws.Cells.Item(479, 15).Value = 16000
$ws.Cells.Item(479, 16).Value = 16000
$ws.Cells.Item(479, 19).Value = 800
$ws.Cells.Item(480, 4).Value = 44617
$ws.Cells.Item(480, 12).Value = "Primera Pintón"
$ws.Cells.Item(480, 13).Value = 240
$ws.Cells.Item(480, 14).Value = 17000
$ws.Cells.Item(480, 15).Value = 18000
$ws.Cells.Item(480, 16).Value = 17500
$ws.Cells.Item(480, 19).Value = 875
$ws.Cells.Item(481, 4).Value = 44557
$ws.Cells.Item(481, 12).Value = "Pintón"
$ws.Cells.Item(481, 13).Value = 80
$ws.Cells.Item(481, 14).Value = 11000
$ws.Cells.Item(481, 15).Value = 11000
$ws.Cells.Item(481, 16).Value = 11000
$ws.Cells.Item(481, 19).Value = 550
$ws.Cells.Item(482, 4).Value = 44557
$ws.Cells.Item(482, 12).Value = "Primera Pintón"
$ws.Cells.Item(482, 13).Value = 240
$ws.Cells.Item(482, 14).Value = 12000
$ws.Cells.Item(482, 15).Value = 13000
$ws.Cells.Item(482, 16).Value = 12500
$ws.Cells.Item(482, 19).Value = 625
$ws.Cells.Item(483, 4).Value = 44489
$ws.Cells.Item(483, 12).Value = "Primera Pintón"
$ws.Cells.Item(483, 13).Value = 400
$ws.Cells.Item(483, 14).Value = 26000
$ws.Cells.Item(483, 15).Value = 27000
$ws.Cells.Item(483, 16).Value = 26500
$ws.Cells.Item(483, 19).Value = 1325
$ws.Cells.Item(484, 4).Value = 44264
$ws.Cells.Item(484, 12).Value = "Pintón"
$ws.Cells.Item(484, 13).Value = 80
$ws.Cells.Item(484, 14).Value = 14000
$ws.Cells.Item(484, 15).Value = 14000
$ws.Cells.Item(484, 16).Value = 14000
$ws.Cells.Item(484, 19).Value = 700
$ws.Cells.Item(485, 4).Value = 44264
$ws.Cells.Item(485, 12).Value = "Primera Pintón"
$ws.Cells.Item(485, 13).Value = 360
$ws.Cells.Item(485, 14).Value = 15000
$ws.Cells.Item(485, 15).Value = 16000
$ws.Cells.Item(485, 16).Value = 15500
$ws.Cells.Item(485, 19).Value = 775
$ws.Cells.Item(486, 4).Value = 44396
$ws.Cells.Item(486, 12).Value = "Pintón"
$ws.Cells.Item(486, 13).Value = 80
$ws.Cells.Item(486, 14).Value = 11000
$ws.Cells.Item(486, 15).Value = 11000
$ws.Cells.Item(486, 16).Value = 11000
$ws.Cells.Item(486, 19).Value = 550
$ws.Cells.Item(487, 4).Value = 44396
$ws.Cells.Item(487, 12).Value = "Primera Pintón"
$ws.Cells.Item(487, 13).Value = 240
$ws.Cells.Item(487, 14).Value = 12000
$ws.Cells.Item(487, 15).Value = 13000
$ws.Cells.Item(487, 16).Value = 12500
$ws.Cells.Item(487, 19).Value = 625
$ws.Cells.Item(488, 4).Value = 44232
$ws.Cells.Item(488, 12).Value = "Pintón"
$ws.Cells.Item(488, 13).Value = 760
$ws.Cells.Item(488, 14).Value = 14500
$ws.Cells.Item(488, 15).Value = 15000
$ws.Cells.Item(488, 16).Value = 14789
$ws.Cells.Item(488, 19).Value = 739
$ws.Cells.Item(489, 4).Value = 44232
$ws.Cells.Item(489, 12).Value = "Primera Pintón"
$ws.Cells.Item(489, 13).Value = 400
$ws.Cells.Item(489, 14).Value = 15500
$ws.Cells.Item(489, 15).Value = 16000
$ws.Cells.Item(489, 16).Value = 15750
$ws.Cells.Item(489, 19).Value = 788
$ws.Cells.Item(490, 4).Value = 44279
$ws.Cells.Item(490, 12).Value = "Pintón"
$ws.Cells.Item(490, 13).Value = 80
$ws.Cells.Item(490, 14).Value = 11000
$ws.Cells.Item(490, 15).Value = 11000
$ws.Cells.Item(490, 16).Value = 11000
$ws.Cells.Item(490, 19).Value = 550
$ws.Cells.Item(491, 4).Value = 44279
$ws.Cells.Item(491, 12).Value = "Primera Pintón"
$ws.Cells.Item(491, 13).Value = 360
$ws.Cells.Item(491, 14).Value = 12000
$ws.Cells.Item(491, 15).Value = 13000
$ws.Cells.Item(491, 16).Value = 12500
$ws.Cells.Item(491, 19).Value = 625
$ws.Cells.Item(492, 4).Value = 44330
$ws.Cells.Item(492, 12).Value = "Pintón"
$ws.Cells.Item(492, 13).Value = 120
$ws.Cells.Item(492, 14).Value = 12000
$ws.Cells.Item(492, 15).Value = 12000
$ws.Cells.Item(492, 16).Value = 12000
$ws.Cells.Item(492, 19).Value = 600
$ws.Cells.Item(493, 4).Value = 44330
$ws.Cells.Item(493, 12).Value = "Primera Pintón"
$ws.Cells.Item(493, 13).Value = 240
$ws.Cells.Item(493, 14).Value = 13000
$ws.Cells.Item(493, 15).Value = 14000
$ws.Cells.Item(493, 16).Value = 13500
$ws.Cells.Item(493, 19).Value = 675
$ws.Cells.Item(494, 4).Value = 44572
$ws.Cells.Item(494, 12).Value = "Pintón"
$ws.Cells.Item(494, 13).Value = 120
$ws.Cells.Item(494, 14).Value = 11000
$ws.Cells.Item(494, 15).Value = 12000
$ws.Cells.Item(494, 16).Value = 11500
$ws.Cells.Item(494, 19).Value = 575
$ws.Cells.Item(495, 4).Value = 44572
$ws.Cells.Item(495, 12).Value = "Primera Pintón"
$ws.Cells.Item(495, 13).Value = 200
$ws.Cells.Item(495, 14).Value = 13000
$ws.Cells.Item(495, 15).Value = 14000
$ws.Cells.Item(495, 16).Value = 13500
$ws.Cells.Item(495, 19).Value = 675
$ws.Cells.Item(496, 4).Value = 44257
$ws.Cells.Item(496, 12).Value = "Pintón"
$ws.Cells.Item(496, 13).Value = 200
$ws.Cells.Item(496, 14).Value = 13000
$ws.Cells.Item(496, 15).Value = 13000
$ws.Cells.Item(496, 16).Value = 13000
$ws.Cells.Item(496, 19).Value = 650
$ws.Cells.Item(497, 4).Value = 44257
$ws.Cells.Item(497, 12).Value = "Primera Pintón"
$ws.Cells.Item(497, 13).Value = 360
$ws.Cells.Item(497, 14).Value = 14000
$ws.Cells.Item(497, 15).Value = 15000
$ws.Cells.Item(497, 16).Value = 14500
$ws.Cells.Item(497, 19).Value = 725
$ws.Cells.Item(498, 4).Value = 44301
$ws.Cells.Item(498, 12).Value = "Pintón"
$ws.Cells.Item(498, 13).Value = 80
$ws.Cells.Item(498, 14).Value = 13000
$ws.Cells.Item(498, 15).Value = 13000
$ws.Cells.Item(498, 16).Value = 13000
$ws.Cells.Item(498, 19).Value = 650
$ws.Cells.Item(499, 4).Value = 44301
$ws.Cells.Item(499, 12).Value = "Primera Pintón"
$ws.Cells.Item(499, 13).Value = 160
$ws.Cells.Item(499, 14).Value = 14000
$ws.Cells.Item(499, 15).Value = 15000
$ws.Cells.Item(499, 16).Value = 14500
$ws.Cells.Item(499, 19).Value = 725
$ws.Cells.Item(500, 4).Value = 44370
$ws.Cells.Item(500, 12).Value = "Pintón"
$ws.Cells.Item(500, 13).Value = 60
$ws.Cells.Item(500, 14).Value = 11000
$ws.Cells.Item(500, 15).Value = 11000
$ws.Cells.Item(500, 16).Value = 11000
$ws.Cells.Item(500, 19).Value = 550
$ws.Cells.Item(501, 4).Value = 44370
$ws.Cells.Item(501, 12).Value = "Primera Pintón"
$ws.Cells.Item(501, 13).Value = 120
$ws.Cells.Item(501, 14).Value = 12000
$ws.Cells.Item(501, 15).Value = 13000
$ws.Cells.Item(501, 16).Value = 12500
$ws.Cells.Item(501, 19).Value = 625
$ws.Cells.Item(502, 4).Value = 44487
$ws.Cells.Item(502, 12).Value = "Pintón"
$ws.Cells.Item(502, 13).Value = 120
$ws.Cells.Item(502, 14).Value = 24000
$ws.Cells.Item(502, 15).Value = 24000
$ws.Cells.Item(502, 16).Value = 24000
$ws.Cells.Item(502, 19).Value = 1200
$ws.Cells.Item(503, 4).Value = 44487
$ws.Cells.Item(503, 12).Value = "Primera Pintón"
$ws.Cells.Item(503, 13).Value = 360
$ws.Cells.Item(503, 14).Value = 25000
$ws.Cells.Item(503, 15).Value = 26000
$ws.Cells.Item(503, 16).Value = 25500
$ws.Cells.Item(503, 19).Value = 1275
$ws.Cells.Item(504, 4).Value = 44174
$ws.Cells.Item(504, 12).Value = "Pintón"
$ws.Cells.Item(504, 13).Value = 640
$ws.Cells.Item(504, 14).Value = 14000
$ws.Cells.Item(504, 15).Value = 15000
$ws.Cells.Item(504, 16).Value = 14469
$ws.Cells.Item(504, 19).Value = 723
$ws.Cells.Item(505, 4).Value = 44174
$ws.Cells.Item(505, 12).Value = "Primera Pintón"
$ws.Cells.Item(505, 13).Value = 300
$ws.Cells.Item(505, 14).Value = 15500
$ws.Cells.Item(505, 15).Value = 16000
$ws.Cells.Item(505, 16).Value = 15667
$ws.Cells.Item(505, 19).Value = 783
$ws.Cells.Item(506, 4).Value = 44200
$ws.Cells.Item(506, 12).Value = "Pintón"
$ws.Cells.Item(506, 13).Value = 160
$ws.Cells.Item(506, 14).Value = 14000
$ws.Cells.Item(506, 15).Value = 15000
$ws.Cells.Item(506, 16).Value = 14500
$ws.Cells.Item(506, 19).Value = 725
$ws.Cells.Item(507, 4).Value = 44200
$ws.Cells.Item(507, 12).Value = "Primera Pintón"
$ws.Cells.Item(507, 13).Value = 240
$ws.Cells.Item(507, 14).Value = 14000
$ws.Cells.Item(507, 15).Value = 15000
$ws.Cells.Item(507, 16).Value = 14500
$ws.Cells.Item(507, 19).Value = 725
$ws.Cells.Item(508, 4).Value = 44385
$ws.Cells.Item(508, 12).Value = "Pintón"
$ws.Cells.Item(508, 13).Value = 80
$ws.Cells.Item(508, 14).Value = 11000
$ws.Cells.Item(508, 15).Value = 11000
$ws.Cells.Item(508, 16).Value = 11000
$ws.Cells.Item(508, 19).Value = 550
$ws.Cells.Item(509, 4).Value = 44385
$ws.Cells.Item(509, 12).Value = "Primera Pintón"
$ws.Cells.Item(509, 13).Value = 240
$ws.Cells.Item(509, 14).Value = 12000
$ws.Cells.Item(509, 15).Value = 13000
$ws.Cells.Item(509, 16).Value = 12500
$ws.Cells.Item(509, 19).Value = 625
$ws.Cells.Item(510, 4).Value = 44236
$ws.Cells.Item(510, 12).Value = "Pintón"
$ws.Cells.Item(510, 13).Value = 80
$ws.Cells.Item(510, 14).Value = 14000
$ws.Cells.Item(510, 15).Value = 14000
$ws.Cells.Item(510, 16).Value = 14000
$ws.Cells.Item(510, 19).Value = 700
$ws.Cells.Item(511, 4).Value = 44236
$ws.Cells.Item(511, 12).Value = "Primera Pintón"
$ws.Cells.Item(511, 13).Value = 240
$ws.Cells.Item(511, 14).Value = 15000
$ws.Cells.Item(511, 15).Value = 16000
$ws.Cells.Item(511, 16).Value = 15500
$ws.Cells.Item(511, 19).Value = 775
$ws.Cells.Item(512, 4).Value = 44221
$ws.Cells.Item(512, 12).Value = "Pintón"
$ws.Cells.Item(512, 13).Value = 80
$ws.Cells.Item(512, 14).Value = 12000
$ws.Cells.Item(512, 15).Value = 12000
$ws.Cells.Item(512, 16).Value = 12000
$ws.Cells.Item(512, 19).Value = 600
$ws.Cells.Item(513, 4).Value = 44221
$ws.Cells.Item(513, 12).Value = "Primera Pintón"
$ws.Cells.Item(513, 13).Value = 240
$ws.Cells.Item(513, 14).Value = 13000
$ws.Cells.Item(513, 15).Value = 14000
$ws.Cells.Item(513, 16).Value = 13500
$ws.Cells.Item(513, 19).Value = 675
$ws.Cells.Item(514, 4).Value = 44413
$ws.Cells.Item(514, 12).Value = "Pintón"
$ws.Cells.Item(514, 13).Value = 200
$ws.Cells.Item(514, 14).Value = 11000
$ws.Cells.Item(514, 15).Value = 11000
$ws.Cells.Item(514, 16).Value = 11000
$ws.Cells.Item(514, 19).Value = 550
$ws.Cells.Item(515, 4).Value = 44413
$ws.Cells.Item(515, 12).Value = "Primera Pintón"
$ws.Cells.Item(515, 13).Value = 600
$ws.Cells.Item(515, 14).Value = 12000
$ws.Cells.Item(515, 15).Value = 13000
$ws.Cells.Item(515, 16).Value = 12500
$ws.Cells.Item(515, 19).Value = 625
$ws.Cells.Item(516, 4).Value = 44229
$ws.Cells.Item(516, 12).Value = "Pintón"
$ws.Cells.Item(516, 13).Value = 590
$ws.Cells.Item(516, 14).Value = 11500
$ws.Cells.Item(516, 15).Value = 12000
$ws.Cells.Item(516, 16).Value = 11703
$ws.Cells.Item(516, 19).Value = 585
$ws.Cells.Item(517, 4).Value = 44229
$ws.Cells.Item(517, 12).Value = "Primera Pintón"
$ws.Cells.Item(517, 13).Value = 400
$ws.Cells.Item(517, 14).Value = 12500
$ws.Cells.Item(517, 15).Value = 13000
$ws.Cells.Item(517, 16).Value = 12812
$ws.Cells.Item(517, 19).Value = 641
$ws.Cells.Item(518, 4).Value = 44214
$ws.Cells.Item(518, 12).Value = "Pintón"
$ws.Cells.Item(518, 13).Value = 740
$ws.Cells.Item(518, 14).Value = 10000
$ws.Cells.Item(518, 15).Value = 10500
$ws.Cells.Item(518, 16).Value = 10230
$ws.Cells.Item(518, 19).Value = 512
$ws.Cells.Item(519, 4).Value = 44214
$ws.Cells.Item(519, 12).Value = "Primera Pintón"
$ws.Cells.Item(519, 13).Value = 350
$ws.Cells.Item(519, 14).Value = 11000
$ws.Cells.Item(519, 15).Value = 12000
$ws.Cells.Item(519, 16).Value = 11429
$ws.Cells.Item(519, 19).Value = 571
$ws.Cells.Item(520, 4).Value = 44299
$ws.Cells.Item(520, 12).Value = "Pintón"
$ws.Cells.Item(520, 13).Value = 80
$ws.Cells.Item(520, 14).Value = 13000
$ws.Cells.Item(520, 15).Value = 13000
$ws.Cells.Item(520, 16).Value = 13000
$ws.Cells.Item(520, 19).Value = 650
$ws.Cells.Item(521, 4).Value = 44299
$ws.Cells.Item(521, 12).Value = "Primera Pintón"
$ws.Cells.Item(521, 13).Value = 160
$ws.Cells.Item(521, 14).Value = 14000
$ws.Cells.Item(521, 15).Value = 15000
$ws.Cells.Item(521, 16).Value = 14500
$ws.Cells.Item(521, 19).Value = 725
$ws.Cells.Item(522, 4).Value = 44610
$ws.Cells.Item(522, 12).Value = "Pintón"
$ws.Cells.Item(522, 13).Value = 500
$ws.Cells.Item(522, 14).Value = 16000
$ws.Cells.Item(522, 15).Value = 16000
$ws.Cells.Item(522, 16).Value = 16000
$ws.Cells.Item(522, 19).Value = 800
$ws.Cells.Item(523, 4).Value = 44610
$ws.Cells.Item(523, 12).Value = "Primera Pintón"
$ws.Cells.Item(523, 13).Value = 600
$ws.Cells.Item(523, 14).Value = 17000
$ws.Cells.Item(523, 15).Value = 18000
$ws.Cells.Item(523, 16).Value = 17500
$ws.Cells.Item(523, 19).Value = 875
$ws.Cells.Item(524, 4).Value = 44312
$ws.Cells.Item(524, 12).Value = "Pintón"
$ws.Cells.Item(524, 13).Value = 120
$ws.Cells.Item(524, 14).Value = 11000
$ws.Cells.Item(524, 15).Value = 11000
$ws.Cells.Item(524, 16).Value = 11000
$ws.Cells.Item(524, 19).Value = 550
$ws.Cells.Item(525, 4).Value = 44312
$ws.Cells.Item(525, 12).Value = "Primera Pintón"
$ws.Cells.Item(525, 13).Value = 240
$ws.Cells.Item(525, 14).Value = 12000
$ws.Cells.Item(525, 15).Value = 13000
$ws.Cells.Item(525, 16).Value = 12500
$ws.Cells.Item(525, 19).Value = 625
$ws.Cells.Item(526, 4).Value = 44399
$ws.Cells.Item(526, 12).Value = "Primera Pintón"
$ws.Cells.Item(526, 13).Value = 240
$ws.Cells.Item(526, 14).Value = 16000
$ws.Cells.Item(526, 15).Value = 17000
$ws.Cells.Item(526, 16).Value = 16500
$ws.Cells.Item(526, 19).Value = 825
$ws.Cells.Item(527, 4).Value = 44615
$ws.Cells.Item(527, 12).Value = "Pintón"
$ws.Cells.Item(527, 13).Value = 80
$ws.Cells.Item(527, 14).Value = 15000
$ws.Cells.Item(527, 15).Value = 15000
$ws.Cells.Item(527, 16).Value = 15000
$ws.Cells.Item(527, 19).Value = 750
$ws.Cells.Item(528, 4).Value = 44615
$ws.Cells.Item(528, 12).Value = "Primera Pintón"
$ws.Cells.Item(528, 13).Value = 240
$ws.Cells.Item(528, 14).Value = 16000
$ws.Cells.Item(528, 15).Value = 17000
$ws.Cells.Item(528, 16).Value = 16500
$ws.Cells.Item(528, 19).Value = 825
$ws.Cells.Item(529, 4).Value = 44522
$ws.Cells.Item(529, 12).Value = "Pintón"
$ws.Cells.Item(529, 13).Value = 60
$ws.Cells.Item(529, 14).Value = 15000
$ws.Cells.Item(529, 15).Value = 15000
$ws.Cells.Item(529, 16).Value = 15000
$ws.Cells.Item(529, 19).Value = 750
$ws.Cells.Item(530, 4).Value = 44522
$ws.Cells.Item(530, 12).Value = "Primera Pintón"
$ws.Cells.Item(530, 13).Value = 160
$ws.Cells.Item(530, 14).Value = 16000
$ws.Cells.Item(530, 15).Value = 17000
$ws.Cells.Item(530, 16).Value = 16500
$ws.Cells.Item(530, 19).Value = 825
$ws.Cells.Item(531, 4).Value = 44543
$ws.Cells.Item(531, 12).Value = "Pintón"
$ws.Cells.Item(531, 13).Value = 120
$ws.Cells.Item(531, 14).Value = 15000
$ws.Cells.Item(531, 15).Value = 15000
$ws.Cells.Item(531, 16).Value = 15000
$ws.Cells.Item(531, 19).Value = 750
$ws.Cells.Item(532, 4).Value = 44543
$ws.Cells.Item(532, 12).Value = "Primera Pintón"
$ws.Cells.Item(532, 13).Value = 360
$ws.Cells.Item(532, 14).Value = 16000
$ws.Cells.Item(532, 15).Value = 17000
$ws.Cells.Item(532, 16).Value = 16500
$ws.Cells.Item(532, 19).Value = 825
$ws.Cells.Item(533, 4).Value = 44167
$ws.Cells.Item(533, 12).Value = "Pintón"
$ws.Cells.Item(533, 13).Value = 80
$ws.Cells.Item(533, 14).Value = 16000
$ws.Cells.Item(533, 15).Value = 16000
$ws.Cells.Item(533, 16).Value = 16000
$ws.Cells.Item(533, 19).Value = 800
$ws.Cells.Item(534, 4).Value = 44167
$ws.Cells.Item(534, 12).Value = "Primera Pintón"
$ws.Cells.Item(534, 13).Value = 160
$ws.Cells.Item(534, 14).Value = 17000
$ws.Cells.Item(534, 15).Value = 18000
$ws.Cells.Item(534, 16).Value = 17500
$ws.Cells.Item(534, 19).Value = 875
$ws.Cells.Item(535, 4).Value = 44277
$ws.Cells.Item(535, 12).Value = "Pintón"
$ws.Cells.Item(535, 13).Value = 120
$ws.Cells.Item(535, 14).Value = 12000
$ws.Cells.Item(535, 15).Value = 12000
$ws.Cells.Item(535, 16).Value = 12000
$ws.Cells.Item(535, 19).Value = 600
$ws.Cells.Item(536, 4).Value = 44277
$ws.Cells.Item(536, 12).Value = "Primera Pintón"
$ws.Cells.Item(536, 13).Value = 240
$ws.Cells.Item(536, 14).Value = 13000
$ws.Cells.Item(536, 15).Value = 14000
$ws.Cells.Item(536, 16).Value = 13500
$ws.Cells.Item(536, 19).Value = 675
$ws.Cells.Item(537, 4).Value = 44258
$ws.Cells.Item(537, 12).Value = "Primera Pintón"
$ws.Cells.Item(537, 13).Value = 160
$ws.Cells.Item(537, 14).Value = 14500
$ws.Cells.Item(537, 15).Value = 15000
$ws.Cells.Item(537, 16).Value = 14750
$ws.Cells.Item(537, 19).Value = 738
$ws.Cells.Item(538, 4).Value = 44390
$ws.Cells.Item(538, 12).Value = "Pintón"
$ws.Cells.Item(538, 13).Value = 80
$ws.Cells.Item(538, 14).Value = 11000
$ws.Cells.Item(538, 15).Value = 11000
$ws.Cells.Item(538, 16).Value = 11000
$ws.Cells.Item(538, 19).Value = 550
$ws.Cells.Item(539, 4).Value = 44390
$ws.Cells.Item(539, 12).Value = "Primera Pintón"
$ws.Cells.Item(539, 13).Value = 240
$ws.Cells.Item(539, 14).Value = 13000
$ws.Cells.Item(539, 15).Value = 14000
$ws.Cells.Item(539, 16).Value = 13500
$ws.Cells.Item(539, 19).Value = 675
$ws.Cells.Item(540, 4).Value = 44349
$ws.Cells.Item(540, 12).Value = "Pintón"
$ws.Cells.Item(540, 13).Value = 120
$ws.Cells.Item(540, 14).Value = 10000
$ws.Cells.Item(540, 15).Value = 11000
$ws.Cells.Item(540, 16).Value = 10500
$ws.Cells.Item(540, 19).Value = 525
$ws.Cells.Item(541, 4).Value = 44349
$ws.Cells.Item(541, 12).Value = "Primera Pintón"
$ws.Cells.Item(541, 13).Value = 360
$ws.Cells.Item(541, 14).Value = 13000
$ws.Cells.Item(541, 15).Value = 14000
$ws.Cells.Item(541, 16).Value = 13500
$ws.Cells.Item(541, 19).Value = 675
$ws.Cells.Item(542, 4).Value = 44285
$ws.Cells.Item(542, 12).Value = "Pintón"
$ws.Cells.Item(542, 13).Value = 120
$ws.Cells.Item(542, 14).Value = 12000
$ws.Cells.Item(542, 15).Value = 12000
$ws.Cells.Item(542, 16).Value = 12000
$ws.Cells.Item(542, 19).Value = 600
$ws.Cells.Item(543, 4).Value = 44285
$ws.Cells.Item(543, 12).Value = "Primera Pintón"
$ws.Cells.Item(543, 13).Value = 360
$ws.Cells.Item(543, 14).Value = 13000
$ws.Cells.Item(543, 15).Value = 14000
$ws.Cells.Item(543, 16).Value = 13500
$ws.Cells.Item(543, 19).Value = 675
$ws.Cells.Item(544, 4).Value = 44498
$ws.Cells.Item(544, 12).Value = "Pintón"
$ws.Cells.Item(544, 13).Value = 160
$ws.Cells.Item(544, 14).Value = 18000
$ws.Cells.Item(544, 15).Value = 19000
$ws.Cells.Item(544, 16).Value = 18500
$ws.Cells.Item(544, 19).Value = 925
$ws.Cells.Item(545, 4).Value = 44498
$ws.Cells.Item(545, 12).Value = "Primera Pintón"
$ws.Cells.Item(545, 13).Value = 240
$ws.Cells.Item(545, 14).Value = 20000
$ws.Cells.Item(545, 15).Value = 21000
$ws.Cells.Item(545, 16).Value = 20500
$ws.Cells.Item(545, 19).Value = 1025
$ws.Cells.Item(546, 4).Value = 44179
$ws.Cells.Item(546, 12).Value = "Pintón"
$ws.Cells.Item(546, 13).Value = 640
$ws.Cells.Item(546, 14).Value = 11500
$ws.Cells.Item(546, 15).Value = 12000
$ws.Cells.Item(546, 16).Value = 11766
$ws.Cells.Item(546, 19).Value = 588
$ws.Cells.Item(547, 4).Value = 44179
$ws.Cells.Item(547, 12).Value = "Primera Pintón"
$ws.Cells.Item(547, 13).Value = 300
$ws.Cells.Item(547, 14).Value = 12500
$ws.Cells.Item(547, 15).Value = 13000
$ws.Cells.Item(547, 16).Value = 12667
$ws.Cells.Item(547, 19).Value = 633
$ws.Cells.Item(548, 4).Value = 44418
$ws.Cells.Item(548, 12).Value = "Pintón"
$ws.Cells.Item(548, 13).Value = 200
$ws.Cells.Item(548, 14).Value = 11000
$ws.Cells.Item(548, 15).Value = 11000
$ws.Cells.Item(548, 16).Value = 11000
$ws.Cells.Item(548, 19).Value = 550
$ws.Cells.Item(549, 4).Value = 44418
$ws.Cells.Item(549, 12).Value = "Primera Pintón"
$ws.Cells.Item(549, 13).Value = 400
$ws.Cells.Item(549, 14).Value = 12000
$ws.Cells.Item(549, 15).Value = 13000
$ws.Cells.Item(549, 16).Value = 12500
$ws.Cells.Item(549, 19).Value = 625
$ws.Cells.Item(550, 4).Value = 44595
$ws.Cells.Item(550, 12).Value = "Pintón"
$ws.Cells.Item(550, 13).Value = 400
$ws.Cells.Item(550, 14).Value = 11000
$ws.Cells.Item(550, 15).Value = 11000
$ws.Cells.Item(550, 16).Value = 11000
$ws.Cells.Item(550, 19).Value = 550
$ws.Cells.Item(551, 4).Value = 44595
$ws.Cells.Item(551, 12).Value = "Primera Pintón"
$ws.Cells.Item(551, 13).Value = 400
$ws.Cells.Item(551, 14).Value = 12000
$ws.Cells.Item(551, 15).Value = 13000
$ws.Cells.Item(551, 16).Value = 12500
$ws.Cells.Item(551, 19).Value = 625
$ws.Cells.Item(552, 4).Value = 44628
$ws.Cells.Item(552, 12).Value = "Pintón"
$ws.Cells.Item(552, 13).Value = 80
$ws.Cells.Item(552, 14).Value = 18000
$ws.Cells.Item(552, 15).Value = 18000
$ws.Cells.Item(552, 16).Value = 18000
$ws.Cells.Item(552, 19).Value = 900
$ws.Cells.Item(553, 4).Value = 44628
$ws.Cells.Item(553, 12).Value = "Primera Pintón"
$ws.Cells.Item(553, 13).Value = 300
$ws.Cells.Item(553, 14).Value = 19000
$ws.Cells.Item(553, 15).Value = 20000
$ws.Cells.Item(553, 16).Value = 19500
$ws.Cells.Item(553, 19).Value = 975
$ws.Cells.Item(554, 4).Value = 44335
$ws.Cells.Item(554, 12).Value = "Pintón"
$ws.Cells.Item(554, 13).Value = 150
$ws.Cells.Item(554, 14).Value = 12000
$ws.Cells.Item(554, 15).Value = 12000
$ws.Cells.Item(554, 16).Value = 12000
$ws.Cells.Item(554, 19).Value = 600
$ws.Cells.Item(555, 4).Value = 44335
$ws.Cells.Item(555, 12).Value = "Primera Pintón"
$ws.Cells.Item(555, 13).Value = 240
$ws.Cells.Item(555, 14).Value = 13000
$ws.Cells.Item(555, 15).Value = 14000
$ws.Cells.Item(555, 16).Value = 13500
$ws.Cells.Item(555, 19).Value = 675
$ws.Cells.Item(556, 4).Value = 44552
$ws.Cells.Item(556, 12).Value = "Pintón"
$ws.Cells.Item(556, 13).Value = 100
$ws.Cells.Item(556, 14).Value = 12000
$ws.Cells.Item(556, 15).Value = 12000
$ws.Cells.Item(556, 16).Value = 12000
$ws.Cells.Item(556, 19).Value = 600
$ws.Cells.Item(557, 4).Value = 44552
$ws.Cells.Item(557, 12).Value = "Primera Pintón"
$ws.Cells.Item(557, 13).Value = 300
$ws.Cells.Item(557, 14).Value = 13000
$ws.Cells.Item(557, 15).Value = 14000
$ws.Cells.Item(557, 16).Value = 13500
$ws.Cells.Item(557, 19).Value = 675
$ws.Cells.Item(558, 4).Value = 44544
$ws.Cells.Item(558, 12).Value = "Pintón"
$ws.Cells.Item(558, 13).Value = 80
$ws.Cells.Item(558, 14).Value = 13000
$ws.Cells.Item(558, 15).Value = 13000
$ws.Cells.Item(558, 16).Value = 13000
$ws.Cells.Item(558, 19).Value = 650
$ws.Cells.Item(559, 4).Value = 44544
$ws.Cells.Item(559, 12).Value = "Primera Pintón"
$ws.Cells.Item(559, 13).Value = 240
$ws.Cells.Item(559, 14).Value = 14000
$ws.Cells.Item(559, 15).Value = 15000
$ws.Cells.Item(559, 16).Value = 14500
$ws.Cells.Item(559, 19).Value = 725
$ws.Cells.Item(560, 4).Value = 44160
$ws.Cells.Item(560, 12).Value = "Pintón"
$ws.Cells.Item(560, 13).Value = 160
$ws.Cells.Item(560, 14).Value = 16000
$ws.Cells.Item(560, 15).Value = 17000
$ws.Cells.Item(560, 16).Value = 16500
$ws.Cells.Item(560, 19).Value = 825
$ws.Cells.Item(561, 4).Value = 44160
$ws.Cells.Item(561, 12).Value = "Primera Pintón"
$ws.Cells.Item(561, 13).Value = 240
$ws.Cells.Item(561, 14).Value = 18000
$ws.Cells.Item(561, 15).Value = 19000
$ws.Cells.Item(561, 16).Value = 18500
$ws.Cells.Item(561, 19).Value = 925
